# chore: update Sheets via scheduled runner
# Refreshes cached market-price-derived figures (currentAveragePrice*,
# Leve cost/profit columns) across the per-job sheets, in place.
$wb = $excel.ActiveWorkbook

$wsALC = $wb.Worksheets.Item("ALC")
$wsALC.Range("H18").Value = 6023.7856
$wsALC.Range("J18").Value = 7819.4
$wsALC.Range("L18").Value = 7819.4
$wsALC.Range("N18").Value = -8387.4
$wsALC.Range("H40").Value = 1761789
$wsALC.Range("I40").Value = 9649.214
$wsALC.Range("K40").Value = 9649.214
$wsALC.Range("M40").Value = -9474.214
$wsALC.Range("H92").Value = 1519.9412
$wsALC.Range("I92").Value = 676.4545000000001
$wsALC.Range("J92").Value = 3066.3333
$wsALC.Range("K92").Value = 676.4545000000001
$wsALC.Range("L92").Value = 3066.3333
$wsALC.Range("M92").Value = 571.5454999999999
$wsALC.Range("N92").Value = -5562.3333
$wsALC.Range("H112").Value = 10474.1
$wsALC.Range("J112").Value = 11379.277
$wsALC.Range("L112").Value = 34137.831
$wsALC.Range("N112").Value = -36353.831
$wsALC.Range("H113").Value = 19237872
$wsALC.Range("I113").Value = 5726.909
$wsALC.Range("K113").Value = 5726.909
$wsALC.Range("M113").Value = -2472.909
$wsALC.Range("H121").Value = 5416.5
$wsALC.Range("J121").Value = 5416.5
$wsALC.Range("L121").Value = 16249.5
$wsALC.Range("N121").Value = -19743.5
$wsALC.Range("H127").Value = 2959.3572
$wsALC.Range("J127").Value = 4683.75
$wsALC.Range("L127").Value = 14051.25
$wsALC.Range("N127").Value = -23971.25
$wsALC.Range("H132").Value = 1436.6522
$wsALC.Range("I132").Value = 1462.75
$wsALC.Range("K132").Value = 4388.25
$wsALC.Range("M132").Value = -1858.25
$wsALC.Range("H138").Value = 3856607.8
$wsALC.Range("I138").Value = 1499.5
$wsALC.Range("J138").Value = 4177866.8
$wsALC.Range("K138").Value = 4498.5
$wsALC.Range("L138").Value = 12533600.4
$wsALC.Range("M138").Value = 641.5
$wsALC.Range("N138").Value = -12543880.4

$wsARM = $wb.Worksheets.Item("ARM")
$wsARM.Range("H32").Value = 3670.3555
$wsARM.Range("I32").Value = 2711.923
$wsARM.Range("K32").Value = 2711.923
$wsARM.Range("M32").Value = -2424.923
$wsARM.Range("H45").Value = 1384
$wsARM.Range("I45").Value = 1180.3529
$wsARM.Range("J45").Value = 1698.7273
$wsARM.Range("K45").Value = 1180.3529
$wsARM.Range("L45").Value = 1698.7273
$wsARM.Range("M45").Value = -803.3529000000001
$wsARM.Range("N45").Value = -2452.7273
$wsARM.Range("H102").Value = 3344.75
$wsARM.Range("I102").Value = 2959.6667
$wsARM.Range("K102").Value = 2959.6667
$wsARM.Range("M102").Value = -1337.6667
$wsARM.Range("H134").Value = 51999.5
$wsARM.Range("J134").Value = 51999.5
$wsARM.Range("L134").Value = 51999.5
$wsARM.Range("N134").Value = -62139.5
$wsARM.Range("H135").Value = 0
$wsARM.Range("J135").Value = 0
$wsARM.Range("L135").Value = 0
$wsARM.Range("N135").ClearContents()

$wsBSM = $wb.Worksheets.Item("BSM")
$wsBSM.Range("H22").Value = 9302.727999999999
$wsBSM.Range("I22").Value = 11350.333
$wsBSM.Range("J22").Value = 88.5
$wsBSM.Range("K22").Value = 11350.333
$wsBSM.Range("L22").Value = 88.5
$wsBSM.Range("M22").Value = -11177.333
$wsBSM.Range("N22").Value = -434.5
$wsBSM.Range("H134").Value = 5628.7188
$wsBSM.Range("I134").Value = 1833.75
$wsBSM.Range("K134").Value = 5501.25
$wsBSM.Range("M134").Value = -2966.25

$wsCRP = $wb.Worksheets.Item("CRP")
$wsCRP.Range("H122").Value = 2446
$wsCRP.Range("I122").Value = 1481
$wsCRP.Range("K122").Value = 4443
$wsCRP.Range("M122").Value = -1993
$wsCRP.Range("H131").Value = 124998.5
$wsCRP.Range("J131").Value = 124998.5
$wsCRP.Range("L131").Value = 124998.5
$wsCRP.Range("N131").Value = -135078.5
$wsCRP.Range("H132").Value = 6208.161
$wsCRP.Range("I132").Value = 4787.125
$wsCRP.Range("K132").Value = 14361.375
$wsCRP.Range("M132").Value = -11831.375
$wsCRP.Range("H134").Value = 6471.5415
$wsCRP.Range("I134").Value = 1594.4286
$wsCRP.Range("K134").Value = 4783.2858
$wsCRP.Range("M134").Value = -2248.2858

$wsCUL = $wb.Worksheets.Item("CUL")
$wsCUL.Range("H107").Value = 1292.2927
$wsCUL.Range("I107").Value = 737.5
$wsCUL.Range("J107").Value = 1820.6666
$wsCUL.Range("K107").Value = 2212.5
$wsCUL.Range("L107").Value = 5461.9998
$wsCUL.Range("M107").Value = -292.5
$wsCUL.Range("N107").Value = -9301.9998
$wsCUL.Range("H132").Value = 4441.278
$wsCUL.Range("I132").Value = 2789.8667
$wsCUL.Range("J132").Value = 12698.333
$wsCUL.Range("K132").Value = 25108.8003
$wsCUL.Range("L132").Value = 114284.997
$wsCUL.Range("M132").Value = -22578.8003
$wsCUL.Range("N132").Value = -119344.997

$wsGSM = $wb.Worksheets.Item("GSM")
$wsGSM.Range("H75").Value = 54999
$wsGSM.Range("J75").Value = 54999
$wsGSM.Range("L75").Value = 54999
$wsGSM.Range("N75").Value = -56747
$wsGSM.Range("H78").Value = 54999
$wsGSM.Range("J78").Value = 54999
$wsGSM.Range("L78").Value = 164997
$wsGSM.Range("N78").Value = -173733

$wsLTW = $wb.Worksheets.Item("LTW")
$wsLTW.Range("H23").Value = 14000
$wsLTW.Range("I23").Value = 14000
$wsLTW.Range("K23").Value = 14000
$wsLTW.Range("M23").Value = -13770
$wsLTW.Range("H46").Value = 3032.5715
$wsLTW.Range("I46").Value = 1521.6
$wsLTW.Range("J46").Value = 4406.1816
$wsLTW.Range("K46").Value = 1521.6
$wsLTW.Range("L46").Value = 4406.1816
$wsLTW.Range("M46").Value = -1333.6
$wsLTW.Range("N46").Value = -4782.1816
$wsLTW.Range("H55").Value = 928.24
$wsLTW.Range("I55").Value = 941
$wsLTW.Range("K55").Value = 941
$wsLTW.Range("M55").Value = -768
$wsLTW.Range("H61").Value = 7142.6665
$wsLTW.Range("I61").Value = 5244.5
$wsLTW.Range("J61").Value = 7685
$wsLTW.Range("K61").Value = 5244.5
$wsLTW.Range("L61").Value = 7685
$wsLTW.Range("M61").Value = -5042.5
$wsLTW.Range("N61").Value = -8089
$wsLTW.Range("H100").Value = 4792.65
$wsLTW.Range("I100").Value = 3591.5833
$wsLTW.Range("K100").Value = 3591.5833
$wsLTW.Range("M100").Value = -3050.5833
$wsLTW.Range("H113").Value = 7142.6665
$wsLTW.Range("I113").Value = 5244.5
$wsLTW.Range("J113").Value = 7685
$wsLTW.Range("K113").Value = 5244.5
$wsLTW.Range("L113").Value = 7685
$wsLTW.Range("M113").Value = -3074.5
$wsLTW.Range("N113").Value = -12025

$wsWVR = $wb.Worksheets.Item("WVR")
$wsWVR.Range("H31").Value = 0
$wsWVR.Range("J31").Value = 0
$wsWVR.Range("L31").Value = 0
$wsWVR.Range("N31").ClearContents()
$wsWVR.Range("H107").Value = 903.8
$wsWVR.Range("I107").Value = 861.8570999999999
$wsWVR.Range("K107").Value = 2585.5713
$wsWVR.Range("M107").Value = -665.5712999999996
$wsWVR.Range("H113").Value = 16693.812
$wsWVR.Range("J113").Value = 1241
$wsWVR.Range("L113").Value = 3723
$wsWVR.Range("N113").Value = -8063
$wsWVR.Range("H114").Value = 56968
$wsWVR.Range("I114").Value = 0
$wsWVR.Range("J114").Value = 56968
$wsWVR.Range("K114").Value = 0
$wsWVR.Range("L114").Value = 56968
$wsWVR.Range("N114").Value = -65646
$wsWVR.Range("M114").ClearContents()
